# Update calculated market/profit figures across multiple crafting-job sheets
# (ALC, ARM, BSM, CRP, GSM, LTW, WVR) per the latest scheduled price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3654.8157
$ws.Range("I64").Value = 3408.7727
$ws.Range("K64").Value = 3408.7727
$ws.Range("M64").Value = -3160.7727

$ws.Range("H67").Value = 3654.8157
$ws.Range("I67").Value = 3408.7727
$ws.Range("K67").Value = 3408.7727
$ws.Range("M67").Value = -2550.7727

$ws.Range("H103").Value = 708.3333
$ws.Range("I103").Value = 720
$ws.Range("J103").Value = 650
$ws.Range("K103").Value = 2160
$ws.Range("L103").Value = 1950
$ws.Range("M103").Value = -1574
$ws.Range("N103").Value = -3122

$ws.Range("H129").Value = 712.66

$ws.Range("H133").Value = 39996.668
$ws.Range("J133").Value = 39996.668
$ws.Range("L133").Value = 39996.668
$ws.Range("N133").Value = -50116.668

$ws.Range("H137").Value = 14634.263
$ws.Range("I137").Value = 809
$ws.Range("J137").Value = 50843.285
$ws.Range("K137").Value = 2427
$ws.Range("L137").Value = 152529.855
$ws.Range("M137").Value = 123
$ws.Range("N137").Value = -157629.855

$ws.Range("H138").Value = 2284.6597
$ws.Range("I138").Value = 1594.8667
$ws.Range("J138").Value = 2415.6328
$ws.Range("K138").Value = 4784.6001
$ws.Range("L138").Value = 7246.8984
$ws.Range("M138").Value = 355.3999000000003
$ws.Range("N138").Value = -17526.8984

$ws.Range("H141").Value = 1378
$ws.Range("I141").Value = 1196.25
$ws.Range("J141").Value = 2105
$ws.Range("K141").Value = 3588.75
$ws.Range("L141").Value = 6315
$ws.Range("M141").Value = 1591.25
$ws.Range("N141").Value = -16675

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10755.756
$ws.Range("I32").Value = 5099.943
$ws.Range("J32").Value = 35499.938
$ws.Range("K32").Value = 5099.943
$ws.Range("L32").Value = 35499.938
$ws.Range("M32").Value = -4812.943
$ws.Range("N32").Value = -36073.938

$ws.Range("H61").Value = 3735.1
$ws.Range("I61").Value = 3285.4285
$ws.Range("J61").Value = 4784.3335
$ws.Range("K61").Value = 3285.4285
$ws.Range("L61").Value = 4784.3335
$ws.Range("M61").Value = -3073.4285
$ws.Range("N61").Value = -5208.3335

$ws.Range("H63").Value = 1683
$ws.Range("I63").Value = 1419.6
$ws.Range("K63").Value = 1419.6
$ws.Range("M63").Value = -733.5999999999999

$ws.Range("H66").Value = 1683
$ws.Range("I66").Value = 1419.6
$ws.Range("K66").Value = 7098
$ws.Range("M66").Value = -3666

$ws.Range("H123").Value = 23966.666
$ws.Range("J123").Value = 23966.666
$ws.Range("L123").Value = 23966.666
$ws.Range("N123").Value = -33766.666

$ws.Range("H128").Value = 32999.668
$ws.Range("J128").Value = 32999.668
$ws.Range("L128").Value = 32999.668
$ws.Range("N128").Value = -42959.668

$ws.Range("H136").Value = 3735.1
$ws.Range("I136").Value = 3285.4285
$ws.Range("J136").Value = 4784.3335
$ws.Range("K136").Value = 9856.2855
$ws.Range("L136").Value = 14353.0005
$ws.Range("M136").Value = -7306.2855
$ws.Range("N136").Value = -19453.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 590.5
$ws.Range("I94").Value = 419.5
$ws.Range("J94").Value = 847
$ws.Range("K94").Value = 419.5
$ws.Range("L94").Value = 847
$ws.Range("M94").Value = 31.5
$ws.Range("N94").Value = -1749

$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws.Range("H134").Value = 3805.6
$ws.Range("I134").Value = 2000
$ws.Range("J134").Value = 5009.3335
$ws.Range("K134").Value = 6000
$ws.Range("L134").Value = 15028.0005
$ws.Range("M134").Value = -3465
$ws.Range("N134").Value = -20098.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1951.9412
$ws.Range("J31").Value = 3533.6667
$ws.Range("L31").Value = 3533.6667
$ws.Range("N31").Value = -4123.6667

$ws.Range("H34").Value = 1951.9412
$ws.Range("J34").Value = 3533.6667
$ws.Range("L34").Value = 3533.6667
$ws.Range("N34").Value = -3937.6667

$ws.Range("H58").Value = 2056.1765
$ws.Range("I58").Value = 1474.4762
$ws.Range("J58").Value = 2995.8462
$ws.Range("K58").Value = 1474.4762
$ws.Range("L58").Value = 2995.8462
$ws.Range("M58").Value = -1271.4762
$ws.Range("N58").Value = -3401.8462

$ws.Range("H132").Value = 3842.2222
$ws.Range("I132").Value = 4266.5557
$ws.Range("J132").Value = 3417.889
$ws.Range("K132").Value = 12799.6671
$ws.Range("L132").Value = 10253.667
$ws.Range("M132").Value = -10269.6671
$ws.Range("N132").Value = -15313.667

$ws.Range("H134").Value = 2432.5789
$ws.Range("I134").Value = 1454.9656
$ws.Range("J134").Value = 3445.1072
$ws.Range("K134").Value = 4364.8968
$ws.Range("L134").Value = 10335.3216
$ws.Range("M134").Value = -1829.8968
$ws.Range("N134").Value = -15405.3216

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws.Range("H136").Value = 2056.1765
$ws.Range("I136").Value = 1474.4762
$ws.Range("J136").Value = 2995.8462
$ws.Range("K136").Value = 4423.4286
$ws.Range("L136").Value = 8987.5386
$ws.Range("M136").Value = -1873.4286
$ws.Range("N136").Value = -14087.5386

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 6331.7393
$ws.Range("I132").Value = 7556.4287
$ws.Range("J132").Value = 4426.6665
$ws.Range("K132").Value = 22669.2861
$ws.Range("L132").Value = 13279.9995
$ws.Range("M132").Value = -20139.2861
$ws.Range("N132").Value = -18339.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3370.818
$ws.Range("I40").Value = 2863.4546
$ws.Range("J40").Value = 3878.182
$ws.Range("K40").Value = 2863.4546
$ws.Range("L40").Value = 3878.182
$ws.Range("M40").Value = -2727.4546
$ws.Range("N40").Value = -4150.182

$ws.Range("H68").Value = 2860.4285
$ws.Range("I68").Value = 2776.2856
$ws.Range("J68").Value = 3112.8572
$ws.Range("K68").Value = 2776.2856
$ws.Range("L68").Value = 3112.8572
$ws.Range("M68").Value = -2027.2856
$ws.Range("N68").Value = -4610.8572

$ws.Range("H71").Value = 2860.4285
$ws.Range("I71").Value = 2776.2856
$ws.Range("J71").Value = 3112.8572
$ws.Range("K71").Value = 13881.428
$ws.Range("L71").Value = 15564.286
$ws.Range("M71").Value = -10137.428
$ws.Range("N71").Value = -23052.286

$ws.Range("H87").Value = 30990
$ws.Range("J87").Value = 30990
$ws.Range("L87").Value = 30990
$ws.Range("N87").Value = -33236

$ws.Range("H90").Value = 30990
$ws.Range("J90").Value = 30990
$ws.Range("L90").Value = 92970
$ws.Range("N90").Value = -104202

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1536.5962
$ws.Range("I132").Value = 880.5161000000001
$ws.Range("J132").Value = 2505.0952
$ws.Range("K132").Value = 2641.5483
$ws.Range("L132").Value = 7515.285600000001
$ws.Range("M132").Value = -111.5483000000004
$ws.Range("N132").Value = -12575.2856

$ws.Range("H136").Value = 6980.364
$ws.Range("I136").Value = 9029.615
$ws.Range("K136").Value = 27088.845
$ws.Range("M136").Value = -24538.845

Write-Output "Updated 183 cells across 7 worksheets (ALC, ARM, BSM, CRP, GSM, LTW, WVR)"
